$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.841.23"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.097.42"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'246.04"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'55.03"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "'59.19"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "'0.925"
$ws.Range("E13").Value = "  +4.79%  "
$ws.Range("D14").Value = "'15.15"
$ws.Range("E14").Value = "  -6.78%  "
$ws.Range("D15").Value = "2.403.51"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("D17").Value = "2.170.63"
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("D18").Value = "36.796.43"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'17.27"
$ws.Range("E19").Value = "  -5.93%  "
$ws.Range("D20").Value = "'72.86"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "'239.28"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").Value = "'9.77"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'167.31"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'20.99"
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "'5.25"
$ws.Range("E31").Value = "  +7.77%  "
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("E35").Value = "  +9.19%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("D38").Value = "'0.0822"
$ws.Range("E38").Value = "  -7.06%  "
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("D40").Value = "'1.17"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").Value = "'4.91"
$ws.Range("E41").Value = "  -6.55%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "'0.0959"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").Value = "'96.76"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  -9.92%  "
$ws.Range("D46").Value = "1.419.51"
$ws.Range("E46").Value = "  +12.35%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.67"
$ws.Range("E47").Value = "  +13.33%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'16.17"
$ws.Range("E48").Value = "  -6.97%  "
$ws.Range("D49").Value = "'2.47"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").Value = "2.292.13"
$ws.Range("E51").Value = "  +2.80%  "
